$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.317.90"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.922.88"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7428"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.30"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.36"
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3141"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06981"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07997"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7740"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.926.92"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.306"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.83"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.314.74"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.28"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.876"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.59"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007865"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.180.30"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.662"
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.428"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.50"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1276"
$ws.Range("E28").Value = "  -4.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.147"
$ws.Range("E29").Value = "  -6.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.358"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.549"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.360"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.083"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.309"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7525"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.773"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01945"
$ws.Range("E38").Value = "  -1.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.791"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "76.17"
$ws.Range("E41").Value = "  -2.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4471"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.951"
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8415"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.689"
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.62"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.907"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.070.37"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.27"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1219"
$ws.Range("E51").Value = "  +5.03%  "
